$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the "License Information" heading paragraph entirely.
# ---------------------------------------------------------------
$p4 = $d.Paragraphs(4)
if ($p4.Range.Text.TrimEnd([char]13) -eq "License Information") {
    $p4.Range.Delete()
}

# ---------------------------------------------------------------
# 2. Clear out the long license paragraph's inner content (leaving
#    only the leading/trailing empty runs), then rebuild it with
#    the new run sequence described by the diff.
# ---------------------------------------------------------------
$licensePara = $d.Paragraphs(4)
$rng = $licensePara.Range
$oldText = "Questions de Traduction (unfoldingWord) (French) is based on: unfoldingWord® Translation Questions, unfoldingWord, 2022, which is licensed under a CC BY-SA 4.0 license."
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------
# 3. Merge the following paragraph ("This PDF version is provided
#    under the same license.") into this one, deleting its text so
#    only a single trailing empty run is left behind.
# ---------------------------------------------------------------
$licensePara = $d.Paragraphs(4)
$nextPara = $d.Paragraphs(5)
$mergeRng = $d.Range($licensePara.Range.End - 1, $nextPara.Range.End)
$mergeRng.Delete()

# ---------------------------------------------------------------
# 4. Insert the new run sequence just before the trailing empty run,
#    as a single text insertion (to avoid the collapsed-range /
#    trailing-empty-run position ambiguity that otherwise drops
#    character formatting on later inserts), then go back and apply
#    bold formatting only to the first ("unfoldingWord® Translation
#    Questions") segment.
# ---------------------------------------------------------------
$boldText = "unfoldingWord® Translation Questions"
$restText = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. " + `
            "unfoldingWord® Translation Questions" + `
            " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from " + `
            "unfoldingWord® Translation Questions" + `
            " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"
$fullText = $boldText + $restText

$lp = $d.Paragraphs(4)
$pos = $lp.Range.End - 1
$insRng = $d.Range($pos, $pos)
$insRng.InsertAfter($fullText)

$boldRng = $d.Range($pos, $pos + $boldText.Length)
$boldRng.Font.Bold = $true
